{"js": "// Fix rincian item placeholders in the table row so they match the\n// document generator's expected field names:\n//   {{no}}           -> {{rincian_no}}\n//   {{uraian}}       -> {{rincian_uraian}}\n//   {{volume}}       -> {{rincian_volume}}\n//   {{satuan}}       -> {{rincian_satuan}}\n//   {{harga_satuan}} -> {{rincian_harga}}\n//   {{jumlah}}       -> {{rincian_jumlah}}\nconst replacements = [\n  [\"{{no}}\", \"{{rincian_no}}\"],\n  [\"{{uraian}}\", \"{{rincian_uraian}}\"],\n  [\"{{volume}}\", \"{{rincian_volume}}\"],\n  [\"{{satuan}}\", \"{{rincian_satuan}}\"],\n  [\"{{harga_satuan}}\", \"{{rincian_harga}}\"],\n  [\"{{jumlah}}\", \"{{rincian_jumlah}}\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fix rincian item placeholders in the table row so they match the\n# document generator's expected field names:\n#   {{no}}           -> {{rincian_no}}\n#   {{uraian}}       -> {{rincian_uraian}}\n#   {{volume}}       -> {{rincian_volume}}\n#   {{satuan}}       -> {{rincian_satuan}}\n#   {{harga_satuan}} -> {{rincian_harga}}\n#   {{jumlah}}       -> {{rincian_jumlah}}\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"{{no}}\", \"{{rincian_no}}\"),\n    @(\"{{uraian}}\", \"{{rincian_uraian}}\"),\n    @(\"{{volume}}\", \"{{rincian_volume}}\"),\n    @(\"{{satuan}}\", \"{{rincian_satuan}}\"),\n    @(\"{{harga_satuan}}\", \"{{rincian_harga}}\"),\n    @(\"{{jumlah}}\", \"{{rincian_jumlah}}\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
